$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date serial value (45186) for every
# data row (rows 2 through 408). The commit bumps that date by two days to
# 45188 for every one of those rows, leaving everything else untouched.
$ws.Range("C2:C408").Value = 45188
